# Auto update Excel log
# Appends new PIR (Bathroom "No Motion"/"Inactive") rows and new Humidity
# (Bathroom "%"/"Active") rows for 2026-01-30, continuing the existing log.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PIR sheet: append rows 98-110
# ---------------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")

$pirTimes = @(
  "16:11:30","16:11:32","16:11:38","16:11:43","16:11:48","16:11:53","16:11:58",
  "16:12:03","16:12:08","16:12:13","16:12:18","16:12:23","16:12:28"
)

$startRow = 98
for ($i = 0; $i -lt $pirTimes.Length; $i++) {
  $r = $startRow + $i

  # Column A holds a date-formatted string ("2026-01-30"). Excel's value
  # parser auto-converts such literals into real date serials, so force the
  # cell to text first to keep it as the literal string like the rest of
  # the column.
  $pir.Cells.Item($r, 1).NumberFormat = "@"
  $pir.Cells.Item($r, 1).Value = "2026-01-30"

  $pir.Cells.Item($r, 2).Value = $pirTimes[$i]
  $pir.Cells.Item($r, 3).Value = "16:00"
  $pir.Cells.Item($r, 4).Value = "Bathroom"
  $pir.Cells.Item($r, 5).Value = "No Motion"
  $pir.Cells.Item($r, 6).Value = "Inactive"
}

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 77-87
# ---------------------------------------------------------------------------
$hum = $wb.Worksheets.Item("Humidity")

$humTimes = @(
  "16:11:30","16:11:33","16:11:38","16:11:43","16:11:48","16:11:58",
  "16:12:03","16:12:08","16:12:13","16:12:18","16:12:28"
)
$humValues = @(
  "87.5%","87.5%","87.5%","87.5%","87.5%","87.5%","87.5%","87.5%","87.5%","87.5%","87.6%"
)

$startRow = 77
for ($i = 0; $i -lt $humTimes.Length; $i++) {
  $r = $startRow + $i

  $hum.Cells.Item($r, 1).NumberFormat = "@"
  $hum.Cells.Item($r, 1).Value = "2026-01-30"

  $hum.Cells.Item($r, 2).Value = $humTimes[$i]
  $hum.Cells.Item($r, 3).Value = "16:00"
  $hum.Cells.Item($r, 4).Value = "Bathroom"

  # Column E holds a percentage-formatted string ("87.5%"). Force text so it
  # is kept as the literal string instead of being parsed into a numeric
  # percentage value.
  $hum.Cells.Item($r, 5).NumberFormat = "@"
  $hum.Cells.Item($r, 5).Value = $humValues[$i]

  $hum.Cells.Item($r, 6).Value = "Active"
}
